$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds text-formatted numbers (e.g. thousand-dot
# separators like "63.864.23" that are not valid Excel numbers). Force the
# cells we touch to keep a Text number format so COM does not reinterpret
# values such as "35.30" or "14.70" as numeric (which would drop trailing
# zeros) and keeps every Price cell the same inline/text flavour as before.

# --- Row 37 / Row 38 swap: PEPE <-> OKB (name, link, price, volume) ---
$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "53.05"
$ws.Range("E37").Value = "  -4.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0755"
$ws.Range("E38").Value = "  -5.61%  "

# --- Price (D) and Volume(1h) (E) updates for all other rows ---
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.864.23"
$ws.Range("E2").Value = "  -4.17%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.110.69"
$ws.Range("E3").Value = "  -5.35%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.98"
$ws.Range("E5").Value = "  -1.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.07"
$ws.Range("E6").Value = "  -9.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.131.41"
$ws.Range("E8").Value = "  -4.70%  "

# Row 9
$ws.Range("E9").Value = "  -4.64%  "

# Row 10
$ws.Range("E10").Value = "  -7.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("E11").Value = "  -9.32%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  -5.62%  "

# Row 13
$ws.Range("E13").Value = "  -8.22%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.30"
$ws.Range("E14").Value = "  -10.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.620.45"
$ws.Range("E15").Value = "  -5.24%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.116"
$ws.Range("E16").Value = "  +1.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.908.55"
$ws.Range("E17").Value = "  -4.12%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.113.68"
$ws.Range("E18").Value = "  -5.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  -7.99%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.43"
$ws.Range("E20").Value = "  -5.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.70"
$ws.Range("E21").Value = "  -5.83%  "

# Row 22
$ws.Range("E22").Value = "  -7.74%  "

# Row 23
$ws.Range("E23").Value = "  -5.40%  "

# Row 24
$ws.Range("E24").Value = "  -8.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.76"
$ws.Range("E25").Value = "  -3.41%  "

# Row 26
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("E27").Value = "  -8.98%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.47"
$ws.Range("E28").Value = "  -8.58%  "

# Row 29
$ws.Range("E29").Value = "  -12.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.114"
$ws.Range("E30").Value = "  -10.14%  "

# Row 31
$ws.Range("E31").Value = "  -4.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.10%  "

# Row 33
$ws.Range("E33").Value = "  -5.88%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.27"
$ws.Range("E34").Value = "  -6.91%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  -3.04%  "

# Row 36
$ws.Range("E36").Value = "  -9.08%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "458.85"
$ws.Range("E39").Value = "  -8.73%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("E40").Value = "  -16.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0393"
$ws.Range("E41").Value = "  -8.39%  "

# Row 42
$ws.Range("E42").Value = "  -9.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.38"
$ws.Range("E43").Value = "  -5.39%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.855.13"
$ws.Range("E44").Value = "  -5.07%  "

# Row 45
$ws.Range("E45").Value = "  -9.73%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.27"
$ws.Range("E46").Value = "  -13.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  -3.21%  "

# Row 48
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.15"
$ws.Range("E49").Value = "  -10.80%  "

# Row 50
$ws.Range("E50").Value = "  -5.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.26"
$ws.Range("E51").Value = "  -2.71%  "

